$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# 1) Strip the surrounding curly braces from the placeholder texts, editing
#    only the run that holds the placeholder so sibling runs/paragraphs
#    (e.g. the trailing empty paragraph in shape 44) are left untouched.
(Get-ShapeById $s 13).TextFrame.TextRange.Runs(1).Text = "titulo"
(Get-ShapeById $s 14).TextFrame.TextRange.Runs(1).Text = "data"
(Get-ShapeById $s 44).TextFrame.TextRange.Runs(1).Text = "resumo"

# 2) Remove the leftover duplicate placeholder shapes (ids 50, 51, 52) that
#    were accidentally left in the deck with the resolved placeholder text.
(Get-ShapeById $s 52).Delete()
(Get-ShapeById $s 51).Delete()
(Get-ShapeById $s 50).Delete()
